$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 (Hydrogen / Non-metallic minerals) no longer has a value - becomes blank
$ws.Range("D3").Value = $null

# Row 7 used to be "Other" -> it is now relabelled "Biogas" with a corrected value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 367.4123375832285

# A new row 8 is added, taking over the "Other" label with a corrected value.
# Copy the formatting of row 7 (border/bold/alignment) down to row 8 first.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A8:D8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 3781.901685110624
